$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update building cost / numbers ---
$ws.Range("C3").Value = 200
$ws.Range("C4").Value = 200
$ws.Range("C5").Value = 200

$ws.Range("C6").Value = 450
$ws.Range("E6").Value = 150
$ws.Range("H6").Value = 20

$ws.Range("C7").Value = 350
$ws.Range("E7").Value = 250
$ws.Range("H7").Value = 20

$ws.Range("C8").Value = 300
$ws.Range("E8").Value = 100

$ws.Range("D9").Value = 1500
$ws.Range("F9").Value = 1500
$ws.Range("H9").Value = 120
$ws.Range("I9").Value = 0.9

$ws.Range("G10").Value = 3000

# --- Add new tooltip / summary cells (order matters for shared string table) ---
$ws.Range("K22").Value = "2-2-2-2-2-1-1"
$ws.Range("N20").Value = "Workers"
$ws.Range("J23").Value = "Total work:"
$ws.Range("J22").Value = "Build pattern:"
$ws.Range("N21").Value = "Time Limit (mins)"
$ws.Range("N22").Value = "Time Limit (secs)"
$ws.Range("N23").Value = "Total work:"

$ws.Range("O20").Value = 10
$ws.Range("O21").Value = 8
$ws.Range("O22").Formula = "=SUM(O21*60)"
$ws.Range("K23").Formula = "=SUM((K3*2)+(K4*2)+(K5*2)+(K6*2)+(K7*2)+(K8)+(K9) + (J10))"
$ws.Range("O23").Formula = "=SUM(O20 * O22)"

# --- Window / selection state ---
$null = $ws.Range("G11").Select()
